$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 12:52"

# Row 11: Iran - refreshed stats
$ws.Range("B11").Value = 89328
$ws.Range("C11").Value = 1134
$ws.Range("D11").Value = 68193
$ws.Range("E11").Value = 15485
$ws.Range("F11").Value = 3096
$ws.Range("G11").Value = 76
$ws.Range("H11").Value = 5650

# Row 18: Suiza - refreshed stats
$ws.Range("B18").Value = 28894
$ws.Range("C18").Value = 217
$ws.Range("E18").Value = 6301

# Rows 37-39: Catar overtakes Emiratos Arabes Unidos and Indonesia in total
# cases, so the three countries reshuffle order. Catar moves into row 37
# with its newly refreshed stats; Emiratos Arabes Unidos and Indonesia
# shift down into rows 38/39 keeping their existing stats.
$ws.Range("A37").Value = "Catar"
$ws.Range("B37").Value = 9358
$ws.Range("C37").Value = 833
$ws.Range("D37").Value = 929
$ws.Range("E37").Value = 8419
$ws.Range("F37").Value = 72
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 10

$ws.Range("A38").Value = "Emiratos Arabes Unidos"
$ws.Range("B38").Value = 9281
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 1760
$ws.Range("E38").Value = 7457
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 64

$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 8607
$ws.Range("C39").Value = 396
$ws.Range("D39").Value = 1042
$ws.Range("E39").Value = 6845
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 31
$ws.Range("H39").Value = 720

# Row 52: Finlandia - refreshed stats
$ws.Range("B52").Value = 4475
$ws.Range("C52").Value = 80
$ws.Range("E52").Value = 1798

# Row 68: Uzbekistan - refreshed stats
$ws.Range("D68").Value = 689
$ws.Range("E68").Value = 1139

# Row 108: Malta - refreshed stats
$ws.Range("B108").Value = 448
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 249
$ws.Range("E108").Value = 195

# Row 137: Brunei - refreshed stats
$ws.Range("D137").Value = 121
$ws.Range("E137").Value = 16

# Row 166: Nepal - refreshed stats
$ws.Range("D166").Value = 12
$ws.Range("E166").Value = 37
